# Correcting Relevance Markers Walker (2018) - Wolters (2018)
# Updates the record_atd ("C") and mirrored average_simulation_TD ("D")
# columns of the time-to-discovery simulation table with corrected values,
# then refreshes the overall record_atd average in C222.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 84
$ws.Range("D3").Value = 84
$ws.Range("C5").Value = 43
$ws.Range("D5").Value = 43
$ws.Range("C7").Value = 196
$ws.Range("D7").Value = 196
$ws.Range("C9").Value = 133
$ws.Range("D9").Value = 133
$ws.Range("C11").Value = 64
$ws.Range("D11").Value = 64
$ws.Range("C13").Value = 32
$ws.Range("D13").Value = 32
$ws.Range("C15").Value = 28
$ws.Range("D15").Value = 28
$ws.Range("C17").Value = 102
$ws.Range("D17").Value = 102
$ws.Range("C19").Value = 93
$ws.Range("D19").Value = 93
$ws.Range("C21").Value = 61
$ws.Range("D21").Value = 61
$ws.Range("C23").Value = 125
$ws.Range("D23").Value = 125
$ws.Range("C25").Value = 65
$ws.Range("D25").Value = 65
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("C29").Value = 56
$ws.Range("D29").Value = 56
$ws.Range("C31").Value = 145
$ws.Range("D31").Value = 145
$ws.Range("C33").Value = 134
$ws.Range("D33").Value = 134
$ws.Range("C35").Value = 77
$ws.Range("D35").Value = 77
$ws.Range("C37").Value = 19
$ws.Range("D37").Value = 19
$ws.Range("C39").Value = 78
$ws.Range("D39").Value = 78
$ws.Range("C41").Value = 326
$ws.Range("D41").Value = 326
$ws.Range("C43").Value = 870
$ws.Range("D43").Value = 870
$ws.Range("C45").Value = 475
$ws.Range("D45").Value = 475
$ws.Range("C47").Value = 81
$ws.Range("D47").Value = 81
$ws.Range("C49").Value = 53
$ws.Range("D49").Value = 53
$ws.Range("C51").Value = 3
$ws.Range("D51").Value = 3
$ws.Range("C53").Value = 97
$ws.Range("D53").Value = 97
$ws.Range("C55").Value = 72
$ws.Range("D55").Value = 72
$ws.Range("C57").Value = 87
$ws.Range("D57").Value = 87
$ws.Range("C59").Value = 171
$ws.Range("D59").Value = 171
$ws.Range("C62").Value = 220
$ws.Range("D62").Value = 220
$ws.Range("C63").Value = 9
$ws.Range("D63").Value = 9
$ws.Range("C65").Value = 1476
$ws.Range("D65").Value = 1476
$ws.Range("C67").Value = 131
$ws.Range("D67").Value = 131
$ws.Range("C69").Value = 10
$ws.Range("D69").Value = 10
$ws.Range("C71").Value = 124
$ws.Range("D71").Value = 124
$ws.Range("C73").Value = 23
$ws.Range("D73").Value = 23
$ws.Range("C75").Value = 21
$ws.Range("D75").Value = 21
$ws.Range("C77").Value = 352
$ws.Range("D77").Value = 352
$ws.Range("C79").Value = 162
$ws.Range("D79").Value = 162
$ws.Range("C81").Value = 249
$ws.Range("D81").Value = 249
$ws.Range("C83").Value = 118
$ws.Range("D83").Value = 118
$ws.Range("C85").Value = 197
$ws.Range("D85").Value = 197
$ws.Range("C87").Value = 496
$ws.Range("D87").Value = 496
$ws.Range("C89").Value = 991
$ws.Range("D89").Value = 991
$ws.Range("C91").Value = 780
$ws.Range("D91").Value = 780
$ws.Range("C93").Value = 5
$ws.Range("D93").Value = 5
$ws.Range("C95").Value = 691
$ws.Range("D95").Value = 691
$ws.Range("C97").Value = 57
$ws.Range("D97").Value = 57
$ws.Range("C99").Value = 36
$ws.Range("D99").Value = 36
$ws.Range("C101").Value = 25
$ws.Range("D101").Value = 25
$ws.Range("C103").Value = 6
$ws.Range("D103").Value = 6
$ws.Range("C105").Value = 132
$ws.Range("D105").Value = 132
$ws.Range("C107").Value = 42
$ws.Range("D107").Value = 42
$ws.Range("C109").Value = 305
$ws.Range("D109").Value = 305
$ws.Range("C111").Value = 1949
$ws.Range("D111").Value = 1949
$ws.Range("C113").Value = 177
$ws.Range("D113").Value = 177
$ws.Range("C115").Value = 60
$ws.Range("D115").Value = 60
$ws.Range("C117").Value = 14
$ws.Range("D117").Value = 14
$ws.Range("C119").Value = 1813
$ws.Range("D119").Value = 1813
$ws.Range("C121").Value = 70
$ws.Range("D121").Value = 70
$ws.Range("C123").Value = 218
$ws.Range("D123").Value = 218
$ws.Range("C125").Value = 1272
$ws.Range("D125").Value = 1272
$ws.Range("C126").Value = 110
$ws.Range("D126").Value = 110
$ws.Range("C128").Value = 263
$ws.Range("D128").Value = 263
$ws.Range("C130").Value = 17
$ws.Range("D130").Value = 17
$ws.Range("C132").Value = 1366
$ws.Range("D132").Value = 1366
$ws.Range("C134").Value = 39
$ws.Range("D134").Value = 39
$ws.Range("C136").Value = 11
$ws.Range("D136").Value = 11
$ws.Range("C138").Value = 49
$ws.Range("D138").Value = 49
$ws.Range("C140").Value = 20
$ws.Range("D140").Value = 20
$ws.Range("C142").Value = 96
$ws.Range("D142").Value = 96
$ws.Range("C144").Value = 58
$ws.Range("D144").Value = 58
$ws.Range("C146").Value = 304
$ws.Range("D146").Value = 304
$ws.Range("C148").Value = 34
$ws.Range("D148").Value = 34
$ws.Range("C150").Value = 45
$ws.Range("D150").Value = 45
$ws.Range("C152").Value = 195
$ws.Range("D152").Value = 195
$ws.Range("C154").Value = 71
$ws.Range("D154").Value = 71
$ws.Range("C156").Value = 69
$ws.Range("D156").Value = 69
$ws.Range("C158").Value = 66
$ws.Range("D158").Value = 66
$ws.Range("C160").Value = 232
$ws.Range("D160").Value = 232
$ws.Range("C162").Value = 27
$ws.Range("D162").Value = 27
$ws.Range("C164").Value = 1369
$ws.Range("D164").Value = 1369
$ws.Range("C166").Value = 373
$ws.Range("D166").Value = 373
$ws.Range("C168").Value = 1371
$ws.Range("D168").Value = 1371
$ws.Range("C170").Value = 609
$ws.Range("D170").Value = 609
$ws.Range("C172").Value = 2116
$ws.Range("D172").Value = 2116
$ws.Range("C174").Value = 509
$ws.Range("D174").Value = 509
$ws.Range("C176").Value = 219
$ws.Range("D176").Value = 219
$ws.Range("C178").Value = 41
$ws.Range("D178").Value = 41
$ws.Range("C180").Value = 26
$ws.Range("D180").Value = 354.5
$ws.Range("C181").Value = 721
$ws.Range("D181").Value = 721
$ws.Range("C183").Value = 4
$ws.Range("D183").Value = 4
$ws.Range("C185").Value = 89
$ws.Range("D185").Value = 89
$ws.Range("C187").Value = 360
$ws.Range("D187").Value = 360
$ws.Range("C189").Value = 35
$ws.Range("D189").Value = 35
$ws.Range("C191").Value = 15
$ws.Range("D191").Value = 15
$ws.Range("C193").Value = 278
$ws.Range("D193").Value = 278
$ws.Range("C195").Value = 13
$ws.Range("D195").Value = 13
$ws.Range("C197").Value = 410
$ws.Range("D197").Value = 410
$ws.Range("C199").Value = 54
$ws.Range("D199").Value = 54
$ws.Range("C201").Value = 44
$ws.Range("D201").Value = 44
$ws.Range("C203").Value = 167
$ws.Range("D203").Value = 167
$ws.Range("C205").Value = 47
$ws.Range("D205").Value = 47
$ws.Range("C207").Value = 62
$ws.Range("D207").Value = 62
$ws.Range("C209").Value = 29
$ws.Range("D209").Value = 29
$ws.Range("C211").Value = 503
$ws.Range("D211").Value = 503
$ws.Range("C213").Value = 771
$ws.Range("D213").Value = 771
$ws.Range("C215").Value = 992
$ws.Range("D215").Value = 992
$ws.Range("C217").Value = 447
$ws.Range("D217").Value = 447
$ws.Range("C219").Value = 52
$ws.Range("D219").Value = 52
$ws.Range("C221").Value = 30
$ws.Range("D221").Value = 30
$ws.Range("C222").Value = 277.7567567567568
